$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.573.42"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.417.12"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.12%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.421.27"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "4.007.29"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "64.517.90"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "3.403.63"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.75%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.09%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "2.886.91"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0315"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.770"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.37%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
